$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New rows 12-21 appended to the Orders sheet, extending the used range
# from A1:L11 to A1:L21. "A" and "F" columns hold numeric-looking values
# that must stay stored as text (matching the rest of the sheet), so we
# force a text number format on just those cells before writing values.
$textCells = @("A12","F12","F13","F14","A15","F15","F16","F17","A18","F18","F19","F20")
foreach ($cellRef in $textCells) {
  $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("A12").Value = "5"
$ws.Range("C12").Value = "152_白荔枝_White Ohara_Rosa rugosa Thunb._20stems"
$ws.Range("F12").Value = "15"

$ws.Range("C13").Value = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
$ws.Range("F13").Value = "22"

$ws.Range("C14").Value = "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
$ws.Range("F14").Value = "5"

$ws.Range("A15").Value = "6"
$ws.Range("C15").Value = "152_白荔枝_White Ohara_Rosa rugosa Thunb._20stems"
$ws.Range("F15").Value = "24"

$ws.Range("C16").Value = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
$ws.Range("F16").Value = "14"

$ws.Range("C17").Value = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F17").Value = "10"

$ws.Range("A18").Value = "7"
$ws.Range("C18").Value = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F18").Value = "30"

$ws.Range("C19").Value = "152_白荔枝_White Ohara_Rosa rugosa Thunb._20stems"
$ws.Range("F19").Value = "8"

$ws.Range("C20").Value = "157_流沙_Quicksand_Rosa rugosa Thunb._20stems"
$ws.Range("F20").Value = "11"

$ws.Range("C21").Value = "184_微光_shimmer_Rosa rugosa Thunb._20stems"

# Summary sheet: append extra digits onto the existing G2 tracking code.
# This is a long numeric-looking string that must remain text (it already
# overflows numeric precision), so force the text format first.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "015196181942320232115225241410308110"
